$d = $word.ActiveDocument

function Set-ParagraphText($para, [string]$newText) {
    # Replace the entire paragraph's visible text (everything except the
    # trailing paragraph mark) with a single run containing $newText.
    $pr = $para.Range
    $start = $pr.Start
    $end = $pr.End - 1   # exclude the paragraph mark itself
    if ($end -gt $start) {
        $old = $d.Range($start, $end)
        $old.Delete()
    }
    $target = $d.Range($start, $start)
    $target.InsertBefore($newText)
}

function Set-TextByStyle([string]$styleName, [string]$newText) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Style.NameLocal -eq $styleName) {
            Set-ParagraphText $p $newText
            return
        }
    }
}

Set-TextByStyle "Title" "Answers: Introduction to simultaneous equations"
Set-TextByStyle "Author" "Ollie Brooke"
Set-TextByStyle "Abstract" "Answers to questions relating to the guide on introduction to simultaneous equations."
